# DataPack Populator: write the Admin Order Schedule Allocation header row
# plus the hidden "orchestrator queue" helper block (K14:P14) that the
# Populator uses to stage its payment-queue connections, and size the
# columns to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:G1) -----------------------------------------------
$headers = @(
    "Schedule Name",
    "Administrator",
    "Total",
    "Allocated Amount",
    "Available Amount",
    "Payment ID",
    "Full Path"
)

for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
    $cell.NumberFormat = "@"
    $cell.Font.Color = 0
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
}

# --- Orchestrator queue staging cells (row 14, cols K:P) ---------------
$plainCols = @(11, 12, 15, 16)
foreach ($col in $plainCols) {
    $cell = $ws.Cells.Item(14, $col)
    $cell.Font.Color = 0
    $cell.Locked = $true
}

$textCols = @(13, 14)
foreach ($col in $textCols) {
    $cell = $ws.Cells.Item(14, $col)
    $cell.Font.Color = 0
    $cell.Locked = $true
    $cell.NumberFormat = "@"
}

# --- Column widths -------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.42
$ws.Columns.Item(2).ColumnWidth = 19.42
$ws.Columns.Item(3).ColumnWidth = 21.25
$ws.Columns.Item(4).ColumnWidth = 27.59
$ws.Columns.Item(5).ColumnWidth = 33.25
$ws.Columns.Item(6).ColumnWidth = 20.59
$ws.Columns.Item(7).ColumnWidth = 22.42
$ws.Columns.Item(11).ColumnWidth = 13.09
$ws.Columns.Item(12).ColumnWidth = 6.76
$ws.Columns.Item(13).ColumnWidth = 16.42
$ws.Columns.Item(14).ColumnWidth = 16.09
$ws.Columns.Item(15).ColumnWidth = 10.25
$ws.Columns.Item(16).ColumnWidth = 7.76

# --- Selection matches the authoring session's last-known cursor -------
$ws.Range("C12").Select()

Write-Host "DataPack Populator: header + queue staging cells written"
